$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz2")
$ws.Activate()

# Column A becomes a copy of the values already present in column D (A1:A132 = D1:D132)
$dValues = $ws.Range("D1:D132").Value2
$ws.Range("A1:A132").Value2 = $dValues

# New helper column G (rows 2-19) with computed values
$gValues = @(
  0.00006,
  0.03699,
  0.01569,
  0.00055,
  0.80749,
  0.00929,
  0.01171,
  0.00548,
  0.01109,
  0.01433,
  0.00148,
  0.00001,
  0.0443,
  0.01855,
  0.0004,
  0.00303,
  0.01937,
  0.0002
)

for ($i = 0; $i -lt $gValues.Length; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 7).Value2 = $gValues[$i]
}

# Restore the view state captured by the sheet after the edits: scrolled to
# row 16 with the newly filled G2:G19 range selected.
$ws.Range("G2:G19").Select()
$ActiveWindow = $excel.ActiveWindow
$ActiveWindow.ScrollRow = 16
